# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" (fund-holding detail) positioned right
#    before the "总计" (grand-total) summary sheet.
# 2. Insert a new top row into the "总计" sheet summarizing the 2022-Q1
#    holdings (count + market value), pushing the existing quarters down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Create the "2022-Q1" sheet and position it right before "总计".
# ---------------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "2022-Q1"
$total = $wb.Worksheets.Item("总计")
$newSheet.Move($total)

# NOTE: after .Move() the old $newSheet COM reference goes stale in this
# host, so re-resolve the sheet by name before writing to it.
$ws = $wb.Worksheets.Item("2022-Q1")

# Header row, matches the other quarterly sheets (基金代码/基金名称/...).
$headers = @("基金代码", "基金名称", "基金规模", "股票总仓位", "仓位占比", "持有市值(亿元)", "仓位排名")
$headerCols = @("B", "C", "D", "E", "F", "G", "H")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Range($headerCols[$i] + "1")
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# Fund holding rows for 2022-Q1.
$data = @(
    @("012582", "交银施罗德品质增长一年持有期混合型证券投资基金A", "51.62", "92.40", "4.96", "2.5604", 8),
    @("010454", "交银施罗德内需增长一年持有期混合", "35.89", "91.74", "5.04", "1.8089", 7),
    @("005583", "易方达港股通红利灵活配置混合", "6.88", "87.72", "3.81", "0.2621", 9),
    @("012583", "交银施罗德品质增长一年持有期混合型证券投资基金C", "2.37", "92.40", "4.96", "0.1176", 8),
    @("001715", "工银瑞信新焦点灵活配置混合A", "1.42", "89.10", "4.90", "0.0696", 9),
    @("001998", "工银瑞信新焦点灵活配置混合C", "0.33", "89.10", "4.90", "0.0162", 9)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $rec = $data[$i]

    $aCell = $ws.Range("A$row")
    $aCell.Value = $i
    $aCell.Font.Bold = $true
    $aCell.Borders.LineStyle = 1
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160

    # B: 基金代码 - force text so leading zeros survive.
    $bCell = $ws.Range("B$row")
    $bCell.NumberFormat = "@"
    $bCell.Value = $rec[0]

    # C: 基金名称 (text).
    $ws.Range("C$row").Value = $rec[1]

    # D-G: numeric-looking figures that are stored as text in every other
    # quarterly sheet of this workbook - force text the same way.
    $dCell = $ws.Range("D$row")
    $dCell.NumberFormat = "@"
    $dCell.Value = $rec[2]

    $eCell = $ws.Range("E$row")
    $eCell.NumberFormat = "@"
    $eCell.Value = $rec[3]

    $fCell = $ws.Range("F$row")
    $fCell.NumberFormat = "@"
    $fCell.Value = $rec[4]

    $gCell = $ws.Range("G$row")
    $gCell.NumberFormat = "@"
    $gCell.Value = $rec[5]

    # H: 仓位排名 (actual number).
    $ws.Range("H$row").Value = $rec[6]
}

# ---------------------------------------------------------------------------
# 2. Insert a new top data row into "总计" for the 2022-Q1 summary.
# ---------------------------------------------------------------------------
$totalWs = $wb.Worksheets.Item("总计")
$totalWs.Range("A2:D2").Insert(-4121)

$totalA2 = $totalWs.Range("A2")
$totalA2.Value = 0
$totalA2.Font.Bold = $true
$totalA2.Borders.LineStyle = 1
$totalA2.HorizontalAlignment = -4108
$totalA2.VerticalAlignment = -4160

$totalWs.Range("B2").Value = "2022-Q1"
$totalWs.Range("C2").Value = 6
$totalWs.Range("D2").Value = 4.83

# The insert shifted the pre-existing quarters down a row but left their
# running index (column A) untouched - renumber it 0..4 top to bottom.
$totalWs.Range("A3").Value = 1
$totalWs.Range("A4").Value = 2
$totalWs.Range("A5").Value = 3
$totalWs.Range("A6").Value = 4
